$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns retain their original text formatting
# (Excel would otherwise coerce plain numeric-looking strings to numbers,
# dropping trailing zeros or using scientific notation).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '35.755.37'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.892.56'
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.28'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.691'
$ws.Range("E6").Value = '  +0.14%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '43.29'
$ws.Range("E8").Value = '  +0.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '57.22'
$ws.Range("E9").Value = '  +8.79%  '
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0757'
$ws.Range("E11").Value = '  +2.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0987'
$ws.Range("E12").Value = '  +1.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.96'
$ws.Range("E13").Value = '  +14.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.790'
$ws.Range("E14").Value = '  +6.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.169.65'
$ws.Range("E15").Value = '  +0.39%  '
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.902.87'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '35.775.12'
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.21'
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0831'
$ws.Range("E20").Value = '  +1.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '246.55'
$ws.Range("E21").Value = '  +0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '13.03'
$ws.Range("E22").Value = '  +1.62%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.69'
$ws.Range("E24").Value = '  +6.06%  '
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.17'
$ws.Range("E26").Value = '  -0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.43'
$ws.Range("E27").Value = '  +0.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.69'
$ws.Range("E28").Value = '  +2.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.45'
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("E30").Value = '  +0.75%  '
$ws.Range("E31").Value = '  +5.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0607'
$ws.Range("E32").Value = '  +5.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.28'
$ws.Range("E33").Value = '  +2.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").Value = '  +11.91%  '
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("E36").Value = '  -14.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.860'
$ws.Range("E37").Value = '  +0.91%  '
$ws.Range("E38").Value = '  +14.18%  '
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  +4.37%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.16'
$ws.Range("E41").Value = '  +2.53%  '
$ws.Range("E42").Value = '  +27.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.99'
$ws.Range("E43").Value = '  +0.12%  '
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.318.79'
$ws.Range("E45").Value = '  +1.92%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0809'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("E48").Value = '  +0.47%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.75'
$ws.Range("E49").Value = '  +0.55%  '
$ws.Range("E50").Value = '  +1.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '42.98'
$ws.Range("E51").Value = '  -0.55%  '
